# Scheduled runner update: refresh market-board derived profit columns (H-N)
# across the Leve profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3030430.2
$ws.Range("I33").Value = 133.34782
$ws.Range("J33").Value = 12987120
$ws.Range("K33").Value = 133.34782
$ws.Range("L33").Value = 12987120
$ws.Range("M33").Value = 95.65217999999999
$ws.Range("N33").Value = -12987578

$ws.Range("H100").Value = 18520686
$ws.Range("I100").Value = 33335134
$ws.Range("J100").Value = 2625
$ws.Range("K100").Value = 33335134
$ws.Range("L100").Value = 2625
$ws.Range("M100").Value = -33334593
$ws.Range("N100").Value = -3707

$ws.Range("H106").Value = 39218296
$ws.Range("I106").Value = 15154341
$ws.Range("J106").Value = 83335544
$ws.Range("K106").Value = 15154341
$ws.Range("L106").Value = 83335544
$ws.Range("M106").Value = -15153710
$ws.Range("N106").Value = -83336806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4862.737
$ws.Range("I32").Value = 4144.3555
$ws.Range("J32").Value = 7556.6665
$ws.Range("K32").Value = 4144.3555
$ws.Range("L32").Value = 7556.6665
$ws.Range("M32").Value = -3857.3555
$ws.Range("N32").Value = -8130.6665

$ws.Range("H74").Value = 1320.1063
$ws.Range("I74").Value = 1069.8276
$ws.Range("J74").Value = 1723.3334
$ws.Range("K74").Value = 1069.8276
$ws.Range("L74").Value = 1723.3334
$ws.Range("M74").Value = -195.8276000000001
$ws.Range("N74").Value = -3471.3334

$ws.Range("H77").Value = 1320.1063
$ws.Range("I77").Value = 1069.8276
$ws.Range("J77").Value = 1723.3334
$ws.Range("K77").Value = 5349.138000000001
$ws.Range("L77").Value = 8616.666999999999
$ws.Range("M77").Value = -981.1380000000008
$ws.Range("N77").Value = -17352.667

$ws.Range("H132").Value = 2413.9167
$ws.Range("I132").Value = 1443.5862
$ws.Range("J132").Value = 3894.9473
$ws.Range("K132").Value = 4330.7586
$ws.Range("L132").Value = 11684.8419
$ws.Range("M132").Value = -1800.7586
$ws.Range("N132").Value = -16744.8419

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 852.36365
$ws.Range("I80").Value = 361.66666
$ws.Range("J80").Value = 1441.2
$ws.Range("K80").Value = 361.66666
$ws.Range("L80").Value = 1441.2
$ws.Range("M80").Value = 636.33334
$ws.Range("N80").Value = -3437.2

$ws.Range("H83").Value = 852.36365
$ws.Range("I83").Value = 361.66666
$ws.Range("J83").Value = 1441.2
$ws.Range("K83").Value = 1808.3333
$ws.Range("L83").Value = 7206
$ws.Range("M83").Value = 3183.6667
$ws.Range("N83").Value = -17190

$ws.Range("H134").Value = 1740.5306
$ws.Range("I134").Value = 1533.4286
$ws.Range("J134").Value = 2016.6666
$ws.Range("K134").Value = 4600.2858
$ws.Range("L134").Value = 6049.9998
$ws.Range("M134").Value = -2065.2858
$ws.Range("N134").Value = -11119.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 47620570
$ws.Range("I105").Value = 66667692
$ws.Range("K105").Value = 66667692
$ws.Range("M105").Value = -66665945

$ws.Range("H132").Value = 1580.3513
$ws.Range("I132").Value = 1304.742
$ws.Range("K132").Value = 3914.226
$ws.Range("M132").Value = -1384.226

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 999
$ws.Range("J16").Value = 999
$ws.Range("L16").Value = 2997
$ws.Range("N16").Value = -3343

$ws.Range("H33").Value = 11144495
$ws.Range("J33").Value = 50059.168
$ws.Range("L33").Value = 300355.008
$ws.Range("N33").Value = -300921.008

$ws.Range("H68").Value = 3487.4038
$ws.Range("I68").Value = 4605.923
$ws.Range("J68").Value = 2368.8845
$ws.Range("K68").Value = 13817.769
$ws.Range("L68").Value = 7106.6535
$ws.Range("M68").Value = -13006.769
$ws.Range("N68").Value = -8728.6535

$ws.Range("H71").Value = 3487.4038
$ws.Range("I71").Value = 4605.923
$ws.Range("J71").Value = 2368.8845
$ws.Range("K71").Value = 41453.307
$ws.Range("L71").Value = 21319.9605
$ws.Range("M71").Value = -37397.307
$ws.Range("N71").Value = -29431.9605

$ws.Range("H88").Value = 4340
$ws.Range("J88").Value = 4340
$ws.Range("L88").Value = 13020
$ws.Range("N88").Value = -13876

$ws.Range("H91").Value = 4340
$ws.Range("J91").Value = 4340
$ws.Range("L91").Value = 13020
$ws.Range("N91").Value = -15984

$ws.Range("H97").Value = 11111522
$ws.Range("I97").Value = 14286071
$ws.Range("J97").Value = 599.5
$ws.Range("K97").Value = 42858213
$ws.Range("L97").Value = 1798.5
$ws.Range("M97").Value = -42857717
$ws.Range("N97").Value = -2790.5

$ws.Range("H104").Value = 7352.6665
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 7352.6665
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 22057.9995
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -27299.9995

$ws.Range("H113").Value = 151976.94
$ws.Range("I113").Value = 453.47726
$ws.Range("J113").Value = 455023.88
$ws.Range("K113").Value = 1360.43178
$ws.Range("L113").Value = 1365071.64
$ws.Range("M113").Value = 809.5682200000001
$ws.Range("N113").Value = -1369411.64

$ws.Range("H131").Value = 14286893
$ws.Range("I131").Value = 7692799
$ws.Range("J131").Value = 15626319
$ws.Range("K131").Value = 23078397
$ws.Range("L131").Value = 46878957
$ws.Range("M131").Value = -23073357
$ws.Range("N131").Value = -46889037

$ws.Range("H132").Value = 1636257.9
$ws.Range("I132").Value = 674.3077
$ws.Range("J132").Value = 2022850.4
$ws.Range("K132").Value = 6068.7693
$ws.Range("L132").Value = 18205653.6
$ws.Range("M132").Value = -3538.7693
$ws.Range("N132").Value = -18210713.6

$ws.Range("H140").Value = 6673.375
$ws.Range("I140").Value = 6673.375
$ws.Range("K140").Value = 20020.125
$ws.Range("M140").Value = -14840.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 90910190
$ws.Range("I113").Value = 125000980
$ws.Range("J113").Value = 1433.3334
$ws.Range("K113").Value = 125000980
$ws.Range("L113").Value = 1433.3334
$ws.Range("M113").Value = -124998810
$ws.Range("N113").Value = -5773.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2666.6667
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -4082

$ws.Range("H132").Value = 12040391
$ws.Range("I132").Value = 17336696
$ws.Range("J132").Value = 3335.7273
$ws.Range("K132").Value = 52010088
$ws.Range("L132").Value = 10007.1819
$ws.Range("M132").Value = -52007558
$ws.Range("N132").Value = -15067.1819

$ws.Range("H133").Value = 30775.334
$ws.Range("J133").Value = 30775.334
$ws.Range("L133").Value = 30775.334
$ws.Range("N133").Value = -35835.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 17333.334
$ws.Range("I42").Value = 5000
$ws.Range("K42").Value = 5000
$ws.Range("M42").Value = -4622

$ws.Range("H81").Value = 1814.2858
$ws.Range("I81").Value = 1700
$ws.Range("J81").Value = 2100
$ws.Range("K81").Value = 3400
$ws.Range("L81").Value = 4200
$ws.Range("M81").Value = -2339
$ws.Range("N81").Value = -6322

$ws.Range("H84").Value = 1814.2858
$ws.Range("I84").Value = 1700
$ws.Range("J84").Value = 2100
$ws.Range("K84").Value = 17000
$ws.Range("L84").Value = 21000
$ws.Range("M84").Value = -11696
$ws.Range("N84").Value = -31608

$ws.Range("H100").Value = 1090.2858
$ws.Range("I100").Value = 1214.8334
$ws.Range("J100").Value = 343
$ws.Range("K100").Value = 2429.6668
$ws.Range("L100").Value = 686
$ws.Range("M100").Value = -1888.6668
$ws.Range("N100").Value = -1768

$ws.Range("H107").Value = 52633676
$ws.Range("I107").Value = 76923580
$ws.Range("J107").Value = 5534.3335
$ws.Range("K107").Value = 230770740
$ws.Range("L107").Value = 16603.0005
$ws.Range("M107").Value = -230768820
$ws.Range("N107").Value = -20443.0005

$ws.Range("H132").Value = 2178.75
$ws.Range("I132").Value = 2358.4736
$ws.Range("J132").Value = 1916.0769
$ws.Range("K132").Value = 7075.4208
$ws.Range("L132").Value = 5748.2307
$ws.Range("M132").Value = -4545.4208
$ws.Range("N132").Value = -10808.2307
